# Daily attendance processing - 2025-10-27 05:48:29
# Updates "Recorded By" lists (email order), attendance stats, and the
# row-63 session (ANATOMY B1 #2) moving from Not Recorded -> Recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write a literal text value that LOOKS like a number/percent
# (e.g. "14.2%") without Excel's auto-detection converting the cell to a
# numeric percent and re-stamping its style. We do this by entering the
# text as a formula that evaluates to the literal string, then collapsing
# the formula to its value with Paste Special (values only) - this keeps
# the cell's original style/format untouched.
function Set-LiteralText {
    param($rangeAddr, [string]$text)
    $escaped = $text.Replace('"', '""')
    $ws.Range($rangeAddr).Formula = '="' + $escaped + '"'
    $ws.Range($rangeAddr).Copy() | Out-Null
    $ws.Range($rangeAddr).PasteSpecial(-4163) | Out-Null
}

# ---------------------------------------------------------------------
# "Recorded By" (column G) email lists re-ordered (same sets, new order)
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G7").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G17").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G22").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G32").Value = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G33").Value = "rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G47").Value = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G48").Value = "rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G62").Value = "hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G77").Value = "hend_mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("G92").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G97").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G107").Value = "nesmadrahim@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G112").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# ---------------------------------------------------------------------
# Class Statistics block (K3:L10)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 1
Set-LiteralText "L9" "14.2%"
Set-LiteralText "L10" "49.1%"

# ---------------------------------------------------------------------
# Per-group breakdown row 19 (Year 2 / B1)
# ---------------------------------------------------------------------
$ws.Range("O19").Value = 3
$ws.Range("P19").Value = 0
Set-LiteralText "R19" "20.0%"
Set-LiteralText "S19" "42.9%"

# ---------------------------------------------------------------------
# Row 63: ANATOMY, Year 2 / B1, session 2 moved from "Not Recorded" to
# "Recorded". Copy the (already-recorded) formatting from row 62 so the
# row picks up the "Recorded" style (s=2) instead of "Not Recorded" (s=9),
# then fill in the real attendance data.
# ---------------------------------------------------------------------
$ws.Range("A62:I62").Copy() | Out-Null
$ws.Range("A63:I63").PasteSpecial(-4122) | Out-Null
$ws.Range("G63").Value = "mennatulla.medhat@med.asu.edu.eg"
$ws.Range("H63").Value = "81/154"
$ws.Range("I63").Value = "Recorded"

$excel.CutCopyMode = 0
